# Apply edits from diff: update rows 2 & 3, add new row 4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.328321666666668
$ws.Range("H2").Value = 27.984965
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.033075
$ws.Range("N2").Value = 0.099225
$ws.Range("O2").Value = 0.01491216139549877
$ws.Range("P2").Value = 0.01491216139549877
$ws.Range("Q2").Value = 0.3085342391250001
$ws.Range("R2").Value = 2.776808152125
$ws.Range("S2").Value = 0.01491216139549877
$ws.Range("T2").Value = 0.01491216139549877

# Row 3
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.328321666666668
$ws.Range("H3").Value = 27.984965
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.003012
$ws.Range("N3").Value = 3.009036
$ws.Range("O3").Value = 0.4522169864133641
$ws.Range("P3").Value = 0.4522169864133641
$ws.Range("Q3").Value = 9.356418571526667
$ws.Range("R3").Value = 84.20776714374001
$ws.Range("S3").Value = 0.4522169864133641
$ws.Range("T3").Value = 0.4522169864133641

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rbp4"
$ws.Range("C4").Value = "Stra6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.328321666666668
$ws.Range("H4").Value = 27.984965
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.181901333333333
$ws.Range("N4").Value = 3.545704
$ws.Range("O4").Value = 0.5328708521911371
$ws.Range("P4").Value = 0.5328708521911372
$ws.Range("Q4").Value = 11.02515581559556
$ws.Range("R4").Value = 99.22640234036
$ws.Range("S4").Value = 0.5328708521911371
$ws.Range("T4").Value = 0.5328708521911372

